# This workbook's sheet1 is a small Fisher-score worked example. This revision
# rebuilds the "patient drug seq" sample data (columns C & D, rows 2-11) so it
# reflects sequences built *after* the patient's first diagnosis day, and adds
# a second verification block (columns H:J, rows 14-25) that re-runs the same
# Fisher-score computation against that rebuilt data as a sanity check.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Updated source data (f1 = column C, f2 = column D) ----------------
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1

$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1

$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0

$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0

# ---- Second verification block (H:J), mirroring the Label/f1/f2 layout --
$ws.Range("H14").Value = "Label"
$ws.Range("I14").Value = "f1"
$ws.Range("J14").Value = "f2"

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1

$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0

$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0

$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 1

$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1

$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 1

$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 0

$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 0

$ws.Range("H25").Value = "Fscore"
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 0.04

# ---- Refresh the active selection / view to the area of interest --------
$ws.Range("I14").Select()

$excel.Calculate()
